$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3:J3").Cut($ws.Range("L3"))
$ws.Range("A3:J3").Clear()
$ws.Range("A18:J18").Cut($ws.Range("L18"))
$ws.Range("A18:J18").Clear()

$ws.Range("A8:J8").Select()
